$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preencher tempo gasto para as tarefas T1 e T2 (testes de solicitacao de beneficios)
$ws.Range("F3").Value = "1h+30min"
$ws.Range("F2").Value = "50min+1h"

# Atualizar seleção ativa (kanban) para F2
$ws.Range("F2").Select()
